$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "voting"
$ws.Range("I1").Value = "racism"
$ws.Range("J1").Value = "blm"
$ws.Range("K1").Value = "lgbt"
$ws.Range("L1").Value = "indigenous"
$ws.Range("M1").Value = "mentalhealth"

$ws.Range("I9").Select() | Out-Null
